$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(2.179726059491259, 2.089191184858919, 2.035691974967676, 2.014412456098057, 2.010910423113955, 2.035402883675147, 2.148073715398311, 2.385782493052545, 2.570925862391732, 2.657496254550892, 2.690620766824054, 2.683471517211842, 2.660214543484074, 2.646013676394091, 2.565315914551547, 2.516414978595094, 2.488509151719313, 2.479098486112377, 2.521597692579178, 2.667036356425399, 2.764085669485212, 2.712104371091698, 2.51925394014296, 2.31965302548889)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @(0.1228732742888781, 0.10634400533948, 0.0962401253910059, 0.09213327816263472, 0.09145195467550593, 0.09618469726586909, 0.1171641256599969, 0.1587011295827949, 0.18951641796167, 0.2036124836003523, 0.2089625199822649, 0.2078097377753636, 0.2040523850908471, 0.2017525102862692, 0.1885968570953196, 0.1805469389478844, 0.1759241264221032, 0.1743601448846732, 0.1814031055102419, 0.2051556729602169, 0.220750693251631, 0.2124204924417654, 0.1810160165043726, 0.1474160792265309)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colD = @(0.2304411993542459, 0.2204273764007922, 0.2142908285511425, 0.2117925400979175, 0.2113778287503436, 0.2142571268241085, 0.2269855772364338, 0.2520696190260878, 0.2706153454233515, 0.2790866267349372, 0.2823001902312114, 0.2816078286331276, 0.2793508909240074, 0.2779692104940352, 0.2700624840181831, 0.26522134003784, 0.2624400706331187, 0.2614989146068041, 0.2657363490516218, 0.2800136488242231, 0.289378078011282, 0.2843768244980254, 0.2655035073956356, 0.245266265536884)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}

$colE = @(0.05807592456849164, 0.05772313933454853, 0.05751285478055657, 0.05742872881554284, 0.05741485368392141, 0.05751171391688281, 0.05795295624357344, 0.05886958666627073, 0.05957602564257769, 0.0599049403145635, 0.06003060721666742, 0.06000349260095206, 0.05991525651616492, 0.05986135529126635, 0.05955468465093361, 0.05936850623127299, 0.0592621304536074, 0.05922623438870822, 0.05938825164818518, 0.0599411431318515, 0.06030899662805567, 0.0601120615007158, 0.05937932268945723, 0.05861594346739807)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(3.901996374871999, 3.702702403506976, 3.581514854234257, 3.532417022318725, 3.524281425495786, 3.580851553921434, 3.833029101149521, 4.337368848030707, 4.714570110744432, 4.887765828894771, 4.953593037336759, 4.939405041870998, 4.893176564512942, 4.864892098538292, 4.703284591627494, 4.604562428493097, 4.547930452495763, 4.528781339271575, 4.615055926583324, 4.90674833965835, 5.098800512015544, 4.996165573358894, 4.610311428477075, 4.19981096321078)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colH = @(0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429, 0.07973214163530429)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $colH[$i]
}

$colL = @(0.2541781385874629, 0.2461069628530481, 0.2413381715008285, 0.2394415374433834, 0.2391294133550304, 0.2413124042139003, 0.2513561983885779, 0.2725515286155655, 0.2890626982867701, 0.2967835988685295, 0.2997379184986499, 0.2991002862566461, 0.2970260372745486, 0.2957594946665836, 0.2885623767261336, 0.2842012165262702, 0.2817125234066111, 0.2808732714517674, 0.2846634245373139, 0.2976344612406194, 0.3062902338723319, 0.3016540215190417, 0.2844544024031421, 0.2666545708322161)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}
